$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values parse as plain numbers need to be
# forced to Text format so Excel keeps them as strings (matching the
# source data, which stores every Price/Volume cell as inline text),
# then the style is reset back to Normal so no stray cell style lingers.

$ws.Range('D2').Value = '34.468.32'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').Value = '1.805.83'
$ws.Range('E3').Value = '  +0.16%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '224.85'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.602'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.87%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '38.63'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +6.44%  '
$ws.Range('E9').Value = '  -3.62%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0668'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.25%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0983'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.06%  '
$ws.Range('D12').Value = '2.068.13'
$ws.Range('E12').Value = '  +0.23%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.07'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.19%  '
$ws.Range('D14').Value = '1.817.77'
$ws.Range('E14').Value = '  +1.05%  '
$ws.Range('E15').Value = '  -2.37%  '
$ws.Range('D16').Value = '34.466.32'
$ws.Range('E16').Value = '  +0.19%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.37'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.49%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.09'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.29%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '240.92'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.78%  '
$ws.Range('E20').Value = '  -2.87%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.19'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.41%  '
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.09'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.79%  '
$ws.Range('E24').Value = '  +0.72%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '170.77'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.66'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.51'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.83%  '
$ws.Range('E28').Value = '  +3.03%  '
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('E30').Value = '  -1.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.76'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.94%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0513'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.88%  '
$ws.Range('E33').Value = '  -4.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.82'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.65%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.638'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.22%  '
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.05'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.35%  '
$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').Value = '1.301.79'
$ws.Range('E37').Value = '  -6.79%  '
$ws.Range('E38').Value = '  -2.04%  '
$ws.Range('E39').Value = '  -4.58%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '82.60'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.08%  '
$ws.Range('E41').Value = '  +0.82%  '
$ws.Range('E42').Value = '  +3.26%  '
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.953'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.60%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.80'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.64%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.99'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.90%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0513'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.79%  '
$ws.Range('D47').Value = '1.968.86'
$ws.Range('E47').Value = '  +0.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.77'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.29%  '
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '102.71'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.59%  '
$ws.Range('E51').Value = '  -5.93%  '
